# Horarios actualizados Linea 141 - 305
# Refresh the three "scrape" sheets with the latest update timestamp and
# append the newly scraped rows for sheet 1 (LP1912) and sheet 2 (LP1912-215).

$wb = $excel.ActiveWorkbook

$nuevaHora = "02:56:23"

# ---- Sheet 1: LP1912 -------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $nuevaHora"
$ws1.Range("A3").Value = "Total filas: 8"

$ws1.Cells.Item(12, 1).Value = $nuevaHora
$ws1.Cells.Item(12, 2).Value = "04:46"
$ws1.Cells.Item(12, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(12, 4).Value = 110
$ws1.Cells.Item(12, 5).Value = "LP1912"

$ws1.Cells.Item(13, 1).Value = $nuevaHora
$ws1.Cells.Item(13, 2).Value = "04:53"
$ws1.Cells.Item(13, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(13, 4).Value = 117
$ws1.Cells.Item(13, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ---------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $nuevaHora"
$ws2.Range("A3").Value = "Total filas: 4"

$ws2.Cells.Item(9, 1).Value = $nuevaHora
$ws2.Cells.Item(9, 2).Value = "04:46"
$ws2.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(9, 4).Value = 110
$ws2.Cells.Item(9, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
